$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 53 ("Supervision Scheduled  Release Date"),
# pushing the existing row 53 ("Release Date/Time") and everything below it
# down by one.
$ws.Rows.Item(53).Insert()

# The blank row inherited formatting from the row above (row 52, the
# section header "Release/Custody Status"); clear the cells we don't need
# and line the formatted ones (B/C/E) up with the sibling rows in this
# section by cloning their look from the row directly below (old row 53,
# now row 54), matching the original sheet's B/C/E-only row layout.
$ws.Cells.Item(53, 1).Clear()
$ws.Cells.Item(53, 4).Clear()

$ws.Range("B54").Copy()
$ws.Range("B53").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C54").Copy()
$ws.Range("C53").PasteSpecial(-4122)
$ws.Range("E54").Copy()
$ws.Range("E53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(53).RowHeight = 45

# New content: "Supervision Scheduled  Release Date"
$ws.Cells.Item(53, 2).Value = "Supervision Scheduled  Release Date"
$ws.Cells.Item(53, 3).Value = "A date set for a subject's release"
$ws.Cells.Item(53, 5).Value = "/br-doc:BookingReport/j:Detention[@structures:id=/br-doc:BookingReport/j:ActivityCaseAssociation/nc:Activity/@structures:ref]/j:SupervisionAugmentation/j:SupervisionReleaseEligibilityDate/nc:Date"
